# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 14:05"

# --- Swap country names caused by re-ordering of the shared-strings table ---
# Fiyi / Curazao swap (rows 197-198)
$ws.Range("A197").Value = "Curazao"
$ws.Range("A198").Value = "Fiyi"

# Santa Lucia / Nueva Caledonia swap (rows 199, 201) - Belice (row 200) unaffected
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("A201").Value = "Santa Lucia"

# San Bartolome / Bonaire, San Eustaquio y Saba swap (rows 215-216)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"

# --- Update statistic numbers ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1706293
$ws.Range("C4").Value = 67
$ws.Range("D4").Value = 464727
$ws.Range("E4").Value = 1141759
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 99807

# Paises Bajos (row 24)
$ws.Range("B24").Value = 45578
$ws.Range("C24").Value = 133
$ws.Range("G24").Value = 26
$ws.Range("H24").Value = 5856

# Dinamarca (row 49)
$ws.Range("B49").Value = 11428
$ws.Range("C49").Value = 41
$ws.Range("D49").Value = 10044
$ws.Range("E49").Value = 821

# Chequia (row 54)
$ws.Range("B54").Value = 9004
$ws.Range("C54").Value = 2
$ws.Range("E54").Value = 2505

# Kazajistan (row 55)
$ws.Range("D55").Value = 4560
$ws.Range("E55").Value = 4374

# Moldavia (row 63)
$ws.Range("D63").Value = 3884
$ws.Range("E63").Value = 2998
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 265

# Uzbekistan (row 77)
$ws.Range("D77").Value = 2624
$ws.Range("E77").Value = 624

# Republica de Macedonia (row 90)
$ws.Range("B90").Value = 2014
$ws.Range("C90").Value = 15
$ws.Range("D90").Value = 1453
$ws.Range("E90").Value = 445
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 116

# Sri Lanka (row 103)
$ws.Range("B103").Value = 1201
$ws.Range("C103").Value = 19
$ws.Range("E103").Value = 479

# Libano (row 106)
$ws.Range("B106").Value = 1134
$ws.Range("C106").Value = 15
$ws.Range("E106").Value = 420

# Curazao / Fiyi data swap (rows 197-198), following the name swap above
$ws.Range("D197").Value = 14
$ws.Range("H197").Value = 1
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0
